$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.342.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.90%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.943.62"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.98%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'481.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +5.60%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'149.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.65%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.621"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.66%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  -0.03%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -2.53%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  +7.28%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.0000356"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +10.93%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'42.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.39%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'10.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.31%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.562.53"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.06%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'14.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.93%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.957.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +2.84%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'  -0.31%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'19.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.18%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'  -2.90%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'68.415.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.99%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'436.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +1.72%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.42"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +4.93%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'14.58"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.17%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'87.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.74%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'10.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.43%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'  +10.35%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'3.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.36%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'38.46"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +2.23%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +7.29%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'716.88"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.83%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'13.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.53%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'0.130"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -4.00%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +3.36%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("B34").Value = "'PEPE"
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").Value = "'https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").Value = "'0.0" + [char]0x2083 + "0912"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +31.55%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("B35").Value = "'InjectiveProtocol"
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").Value = "'42.23"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.61%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'58.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.85%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -6.78%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.03%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.20%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'2.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +8.40%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  -1.50%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  +11.72%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.56%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.346"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -3.71%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +0.70%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.06%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.49%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.79%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.21"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.88%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'146.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.14%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.11%  "
$ws.Range("E51").Style = "Normal"
